$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1971
$ws.Range("I40").Value = 1919.8
$ws.Range("K40").Value = 1919.8
$ws.Range("M40").Value = -1744.8

$ws.Range("H51").Value = 2670.3333
$ws.Range("I51").Value = 2670.3333
$ws.Range("K51").Value = 2670.3333
$ws.Range("M51").Value = -2186.3333

$ws.Range("H103").Value = 1899.5555
$ws.Range("I103").Value = 1499.4286
$ws.Range("K103").Value = 4498.2858
$ws.Range("M103").Value = -3912.2858

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 3633.6667
$ws.Range("I6").Value = 3633.6667
$ws.Range("K6").Value = 3633.6667
$ws.Range("M6").Value = -3460.6667

$ws.Range("H74").Value = 2374.9167
$ws.Range("I74").Value = 2374.9167
$ws.Range("K74").Value = 2374.9167
$ws.Range("M74").Value = -1500.9167

$ws.Range("H77").Value = 2374.9167
$ws.Range("I77").Value = 2374.9167
$ws.Range("K77").Value = 11874.5835
$ws.Range("M77").Value = -7506.583500000001

$ws.Range("I132").Value = 146803.42
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 440410.26
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -437880.26
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 58500
$ws.Range("J81").Value = 58500
$ws.Range("L81").Value = 58500
$ws.Range("N81").Value = -60622

$ws.Range("H84").Value = 58500
$ws.Range("J84").Value = 58500
$ws.Range("L84").Value = 175500
$ws.Range("N84").Value = -186108

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 1783.875
$ws.Range("I2").Value = 1154.4
$ws.Range("J2").Value = 2833
$ws.Range("K2").Value = 1154.4
$ws.Range("L2").Value = 2833
$ws.Range("M2").Value = -1041.4
$ws.Range("N2").Value = -3059

$ws.Range("H9").Value = 229998
$ws.Range("J9").Value = 229998
$ws.Range("L9").Value = 229998
$ws.Range("N9").Value = -230334

$ws.Range("H69").Value = 25374.5
$ws.Range("I69").Value = 8750
$ws.Range("J69").Value = 41999
$ws.Range("K69").Value = 8750
$ws.Range("L69").Value = 41999
$ws.Range("M69").Value = -8001
$ws.Range("N69").Value = -43497

$ws.Range("H72").Value = 25374.5
$ws.Range("I72").Value = 8750
$ws.Range("J72").Value = 41999
$ws.Range("K72").Value = 26250
$ws.Range("L72").Value = 125997
$ws.Range("M72").Value = -22506
$ws.Range("N72").Value = -133485

$ws.Range("H105").Value = 380
$ws.Range("I105").Value = 380
$ws.Range("K105").Value = 380
$ws.Range("M105").Value = 1367

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1343.3334
$ws.Range("I131").Value = 1343.3334
$ws.Range("J131").Value = 0
$ws.Range("K131").Value = 4030.0002
$ws.Range("L131").Value = 0
$ws.Range("M131").Value = 1009.9998
$ws.Range("N131").ClearContents()

$ws.Range("H137").Value = 3200
$ws.Range("J137").Value = 3200
$ws.Range("L137").Value = 9600
$ws.Range("N137").Value = -19800

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 10005000
$ws.Range("I7").Value = 20000000
$ws.Range("J7").Value = 9999
$ws.Range("K7").Value = 20000000
$ws.Range("L7").Value = 9999
$ws.Range("M7").Value = -19999888
$ws.Range("N7").Value = -10223

$ws.Range("H8").Value = 10005000
$ws.Range("I8").Value = 20000000
$ws.Range("J8").Value = 9999
$ws.Range("K8").Value = 20000000
$ws.Range("L8").Value = 9999
$ws.Range("M8").Value = -19999861
$ws.Range("N8").Value = -10277

$ws.Range("H11").Value = 667999.5
$ws.Range("I11").Value = 667999.5
$ws.Range("K11").Value = 667999.5
$ws.Range("M11").Value = -667860.5

$ws.Range("H132").Value = 6025
$ws.Range("I132").Value = 5733
$ws.Range("K132").Value = 17199
$ws.Range("M132").Value = -14669

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 105
$ws.Range("I2").Value = 105
$ws.Range("K2").Value = 105
$ws.Range("M2").Value = 7

$ws.Range("H22").Value = 2000
$ws.Range("J22").Value = 2000
$ws.Range("L22").Value = 2000
$ws.Range("N22").Value = -2590

$ws.Range("H27").Value = 2000
$ws.Range("J27").Value = 2000
$ws.Range("L27").Value = 2000
$ws.Range("N27").Value = -2214

$ws.Range("H40").Value = 511039.7
$ws.Range("I40").Value = 17374.5
$ws.Range("K40").Value = 17374.5
$ws.Range("M40").Value = -17238.5

$ws.Range("H46").Value = 700
$ws.Range("I46").Value = 700
$ws.Range("K46").Value = 700
$ws.Range("M46").Value = -512

$ws.Range("H48").Value = 0
$ws.Range("I48").Value = 0
$ws.Range("K48").Value = 0
$ws.Range("M48").ClearContents()

$ws.Range("H53").Value = 13034
$ws.Range("I53").Value = 13000
$ws.Range("J53").Value = 13051
$ws.Range("K53").Value = 13000
$ws.Range("L53").Value = 13051
$ws.Range("M53").Value = -12482
$ws.Range("N53").Value = -14087

$ws.Range("H68").Value = 3285.4285
$ws.Range("I68").Value = 2559.8
$ws.Range("J68").Value = 5099.5
$ws.Range("K68").Value = 2559.8
$ws.Range("L68").Value = 5099.5
$ws.Range("M68").Value = -1810.8
$ws.Range("N68").Value = -6597.5

$ws.Range("H71").Value = 3285.4285
$ws.Range("I71").Value = 2559.8
$ws.Range("J71").Value = 5099.5
$ws.Range("K71").Value = 12799
$ws.Range("L71").Value = 25497.5
$ws.Range("M71").Value = -9055
$ws.Range("N71").Value = -32985.5

$ws.Range("H93").Value = 1833.1666
$ws.Range("J93").Value = 1849.75
$ws.Range("L93").Value = 1849.75
$ws.Range("N93").Value = -4345.75

$ws.Range("H122").Value = 3486.5334
$ws.Range("J122").Value = 3688.7778
$ws.Range("L122").Value = 11066.3334
$ws.Range("N122").Value = -15966.3334

$ws.Range("H134").Value = 100500
$ws.Range("J134").Value = 100500
$ws.Range("L134").Value = 100500
$ws.Range("N134").Value = -110640

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 335999.66
$ws.Range("J2").Value = 3999.5
$ws.Range("L2").Value = 3999.5
$ws.Range("N2").Value = -4223.5

$ws.Range("H5").Value = 4472000
$ws.Range("J5").Value = 786666.7
$ws.Range("L5").Value = 786666.7
$ws.Range("N5").Value = -786890.7

$ws.Range("H51").Value = 19549.5
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 19549.5
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 19549.5
$ws.Range("M51").ClearContents()
$ws.Range("N51").Value = -20569.5

$ws.Range("H53").Value = 0
$ws.Range("I53").Value = 0
$ws.Range("K53").Value = 0
$ws.Range("M53").ClearContents()

$ws.Range("H107").Value = 1165.8
$ws.Range("J107").Value = 1990.3334
$ws.Range("L107").Value = 5971.0002
$ws.Range("N107").Value = -9811.0002

$ws.Range("J126").Value = 1000
$ws.Range("L126").Value = 3000
$ws.Range("N126").Value = -7940
